$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment ID_Produit (column O) by 1 for every data row (2..629)
for ($r = 2; $r -le 629; $r++) {
    $ws.Cells.Item($r, 15).Value2 = $ws.Cells.Item($r, 15).Value2 + 1
}

# Row 13: quantity / amount correction (Kamlac évaporé 48x160g sale)
$ws.Cells.Item(13, 10).Value2 = 0.5
$ws.Cells.Item(13, 11).Value2 = 5250.0

# Row 626: quantity / amount correction (Chocolat Jaune stock)
$ws.Cells.Item(626, 10).Value2 = 11.333333
$ws.Cells.Item(626, 11).Value2 = 113333.336

# New row 630: additional stock entry, formatted like the row above it
$ws.Range("A629:P629").Copy()
$ws.Range("A630:P630").PasteSpecial(-4122)

$ws.Cells.Item(630, 1).Value2 = 46027
$ws.Cells.Item(630, 2).Value2 = 0
$ws.Cells.Item(630, 3).Value = "TATA 2"
$ws.Cells.Item(630, 4).Value = "ALIOUNE BADARA (MANSOUR) SANE"
$ws.Cells.Item(630, 5).Value = "Stock Lundi"
$ws.Cells.Item(630, 9).Value = "Kamlac évaporé 48x160g"
$ws.Cells.Item(630, 10).Value2 = 7.0
$ws.Cells.Item(630, 11).Value2 = 73500.0
$ws.Cells.Item(630, 15).Value2 = 134.0
$ws.Cells.Item(630, 16).Value = "S02"
